# "export version 1.3 for windows 4.3.1"
#
# Adds two new header columns (ID / "شناسه" and Group / "گروه") to the
# ware-list sheet, applies number formats to the quantity/price columns,
# moves the selection to the newly inserted column, and switches the page
# to portrait orientation for export/printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (pushes the old H "توضیحات" -> I, and makes room
# for the two new trailing columns of data).
$ws.Columns("H:H").Insert()

# Quantity column (D) gets a 2-decimal numeric format, the price/new
# columns (F:H) get a plain integer numeric format.
$ws.Columns("D:D").NumberFormat = "0.00"
$ws.Columns("F:H").NumberFormat = "0"

# Match the new column's width to its neighbour (G) so the pair reads as
# one visually consistent block.
$ws.Columns("H:H").ColumnWidth = $ws.Columns("G:G").ColumnWidth

# Populate the new headers. "شناسه" (ID) is written first so it lands in
# the shared-string table ahead of "گروه" (Group), matching the column
# order the workbook ends up with (H=گروه, I=توضیحات, J=شناسه).
$ws.Range("J1").Value = "شناسه"
$ws.Range("H1").Value = "گروه"

# Selection moves to the newly-inserted header cell.
[void]$ws.Range("H1").Select()

# Page setup for the Windows export: portrait orientation.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 0

# Best-effort: restore/update the workbook window geometry recorded for
# this export.
$win = $excel.Windows.Item(1)
$win.Left = 2983
$win.Top = 874
$win.Width = 15291
$win.Height = 17640
